$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the 2023 column (O) mirroring the existing year columns ---

# Row 2 (bottom-border spacer row): copy formatting from N2 into O2
$ws.Range("N2").Copy()
$ws.Range("O2").PasteSpecial(-4122)

# Row 3 (year header row): copy formatting from N3 into O3, then set 2023
$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)
$ws.Range("O3").Value = 2023

# Rows 4-6 (data rows): copy formatting from N4:N6 into O4:O6, then set values
$ws.Range("N4:N6").Copy()
$ws.Range("O4:O6").PasteSpecial(-4122)
$ws.Range("O4").Value = 5571
$ws.Range("O5").Value = 74710
$ws.Range("O6").Value = 375715

$excel.CutCopyMode = $false

# --- Row height tweaks to accommodate the extra column ---
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 46.5

# --- Footnote row (7): shrink the font so the note still fits ---
$ws.Range("A7:C7").Font.Size = 8

# --- Capitalize the English footnote text ---
$ws.Range("C7").Value = "*According to the Service for the Regulation and Supervision of the Communications Sector under the Ministry of Digital Development of the Kyrgyz Republic"
